# Apply updated crypto price/volume data per commit "Updated cryptos list on Fri Jul 19 13:46:49 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.711.24'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").Value = '3.422.34'
$ws.Range("E3").Value = '  -1.44%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = '''575.72'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").Value = '''165.08'
$ws.Range("E6").Value = '  +2.72%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").Value = '3.423.67'
$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("D9").Value = '''0.554'
$ws.Range("E9").Value = '  -4.24%  '

$ws.Range("D10").Value = '''7.29'
$ws.Range("E10").Value = '  +1.04%  '

$ws.Range("E11").Value = '  -1.71%  '

$ws.Range("D12").Value = '''0.423'
$ws.Range("E12").Value = '  -3.18%  '

$ws.Range("D13").Value = '4.012.80'
$ws.Range("E13").Value = '  -1.70%  '

$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("D15").Value = '''27.03'
$ws.Range("E15").Value = '  -1.91%  '

$ws.Range("E16").Value = '  -1.63%  '

$ws.Range("D17").Value = '64.707.68'
$ws.Range("E17").Value = '  -0.36%  '

$ws.Range("D18").Value = '3.396.10'
$ws.Range("E18").Value = '  -4.06%  '

$ws.Range("D19").Value = '''6.18'
$ws.Range("E19").Value = '  -0.60%  '

$ws.Range("D20").Value = '''13.56'
$ws.Range("E20").Value = '  -1.99%  '

$ws.Range("D21").Value = '''379.52'
$ws.Range("E21").Value = '  -0.55%  '

$ws.Range("D22").Value = '''7.85'
$ws.Range("E22").Value = '  -1.56%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").Value = '''70.92'
$ws.Range("E24").Value = '  -2.54%  '

$ws.Range("D25").Value = '''0.517'
$ws.Range("E25").Value = '  -2.88%  '

$ws.Range("E26").Value = '  -3.56%  '

$ws.Range("D27").Value = '''9.75'
$ws.Range("E27").Value = '  -1.27%  '

$ws.Range("D28").Value = '''0.179'
$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("E29").Value = '  -0.58%  '

$ws.Range("D30").Value = '''6.13'
$ws.Range("E30").Value = '  +0.16%  '

$ws.Range("D31").Value = '''1.41'
$ws.Range("E31").Value = '  -2.25%  '

$ws.Range("E32").Value = '  -0.28%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''22.97'
$ws.Range("E33").Value = '  -1.84%  '

$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '''0.999'
$ws.Range("E34").Value = '  +0.07%  '

$ws.Range("D35").Value = '''7.06'
$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").Value = '''1.49'
$ws.Range("E36").Value = '  -6.41%  '

$ws.Range("D37").Value = '''159.57'
$ws.Range("E37").Value = '  -0.99%  '

$ws.Range("D38").Value = '''0.870'
$ws.Range("E38").Value = '  +6.96%  '

$ws.Range("D39").Value = '''1.84'
$ws.Range("E39").Value = '  -2.00%  '

$ws.Range("D40").Value = '''0.0725'
$ws.Range("E40").Value = '  -3.67%  '

$ws.Range("D41").Value = '''25.87'
$ws.Range("E41").Value = '  -3.51%  '

$ws.Range("D42").Value = '2.772.84'
$ws.Range("E42").Value = '  -3.19%  '

$ws.Range("D43").Value = '''42.79'
$ws.Range("E43").Value = '  -0.28%  '

$ws.Range("D44").Value = '''6.49'
$ws.Range("E44").Value = '  -0.45%  '

$ws.Range("D45").Value = '''25.85'
$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").Value = '''4.38'
$ws.Range("E46").Value = '  -3.04%  '

$ws.Range("D47").Value = '''0.0305'
$ws.Range("E47").Value = '  -1.46%  '

$ws.Range("D48").Value = '''2.46'
$ws.Range("E48").Value = '  +0.69%  '

$ws.Range("D49").Value = '''330.32'
$ws.Range("E49").Value = '  -0.26%  '

$ws.Range("D50").Value = '''1.05'
$ws.Range("E50").Value = '  -1.30%  '

$ws.Range("D51").Value = '''6.31'
$ws.Range("E51").Value = '  -2.46%  '
